$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for this market/product. It belongs
# chronologically right after the current row 91, so insert a fresh row at
# 92 (pushing the former rows 92:112 down to 93:113) and fill it in.
$ws.Rows.Item(92).Insert()

# Populate the new row 92 with the new record
$ws.Cells.Item(92, 1).Value = 10
$ws.Cells.Item(92, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(92, 3).Value = 'La Araucanía'
$ws.Cells.Item(92, 4).Value = 44785
$ws.Cells.Item(92, 5).Value = 9
$ws.Cells.Item(92, 6).Value = 100112035
$ws.Cells.Item(92, 7).Value = 'Bruselas (repollito)'
$ws.Cells.Item(92, 8).Value = 'Sin especificar'
$ws.Cells.Item(92, 9).Value = 'Primera'
$ws.Cells.Item(92, 10).Value = 50
$ws.Cells.Item(92, 11).Value = 25000
$ws.Cells.Item(92, 12).Value = 25000
$ws.Cells.Item(92, 13).Value = 25000
$ws.Cells.Item(92, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(92, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(92, 16).Value = 2500
$ws.Cells.Item(92, 17).Value = 10
$ws.Cells.Item(92, 18).Value = 'Hortaliza'
